# Improve Excel report formatting: wrap/align the "Vendor/Shop Name" (B) and
# "HSN Codes" (F) columns to top-left alignment, widen those columns, grow the
# row heights to fit the now-multi-line text, and refresh a few data values
# that were re-derived (shortened/deduped HSN code lists, fixed tax total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# XlHAlign / XlVAlign / XlPasteType constants (no enum available in this
# host, so use the literal values Excel itself uses under the hood).
$xlLeft = -4131
$xlTop  = -4160
$xlPasteValues = -4163

# --- Column widths: B 30 -> 35, F 40 -> 45 -------------------------------
# Excel's ColumnWidth property is expressed in characters of the Normal
# style's font and reads ~0.8333 narrower than the width actually stored in
# the OOXML <col> element for this workbook's default font; back that offset
# out so the saved width lands exactly on the target values.
$widthOffset = 5/6
$ws.Columns.Item(2).ColumnWidth = 35 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 45 - $widthOffset

# --- Row heights: rows 2-4 30 -> 45 (to fit the wrapped, multi-line text) ---
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45

# --- Alignment: Vendor/Shop Name (B) and HSN Codes (F) columns -----------
# Was vertical-center + wrap; now left/top + wrap, applied per-column across
# every data row so the style collapses back onto one shared cell format.
$colB = $ws.Range("B2:B4")
$colB.HorizontalAlignment = $xlLeft
$colB.VerticalAlignment = $xlTop
$colB.WrapText = $true

$colF = $ws.Range("F2:F4")
$colF.HorizontalAlignment = $xlLeft
$colF.VerticalAlignment = $xlTop
$colF.WrapText = $true

# --- Updated cell text/values --------------------------------------------
$ws.Range("B2").Value = "WESTSIDE`nSjr Zion, Survey"
$ws.Range("F2").Value = "996211, 62052000, 62052000`n62046200, 48194000, 33072000`n39264099"

$ws.Range("B3").Value = "SONOVISION`nELECTRONICS PVT LTD"

# F3 ("85287219") and J4 ("1506.70") are purely numeric-looking text; a
# plain .Value assignment would be silently coerced into a number (losing
# the trailing zero on J4, and the original's text typing on both). Route
# them through a text formula and then "paste values" over themselves so
# the result is committed as a literal string cell, same as the rest of
# this sheet's text-typed numeric columns, without leaving a residual
# formula or bumping the cell's style (e.g. no quote-prefix flag).
$f3 = $ws.Range("F3")
$f3.Formula = '="85287219"'
$f3.Copy($f3)
$f3.PasteSpecial($xlPasteValues)

$ws.Range("B4").Value = "LAKSHMI`nAGENCIES"
$ws.Range("F4").Value = "15121910, 15121910, 15121910`n15180039, 15180039"

$j4 = $ws.Range("J4")
$j4.Formula = '="1506.70"'
$j4.Copy($j4)
$j4.PasteSpecial($xlPasteValues)
